# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement detail table (rows 16-51) previously listed every
# period for worker CC/7919949 (SAMIR SAER DIAZ) first, followed by every
# period for worker CE/622125973 (JENNIFER ROSEMARY BRITISH CITIZEN).
# The refreshed export groups the two workers by period instead, walking
# the periods from 1607 up to 1712 and, for each period, emitting the
# SAMIR SAER DIAZ row followed by the JENNIFER ROSEMARY BRITISH CITIZEN
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712")

$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "7919949"
    $ws.Cells.Item($row, 4).Value = "SAMIR SAER DIAZ"
    $ws.Cells.Item($row, 5).Value = $p
    $ws.Cells.Item($row, 6).Value = 24640
    $ws.Cells.Item($row, 7).Value = 781242
    $row = $row + 1

    $ws.Cells.Item($row, 2).Value = "CE"
    $ws.Cells.Item($row, 3).Value = "622125973"
    $ws.Cells.Item($row, 4).Value = "JENNIFER ROSEMARY BRITISH CITIZEN"
    $ws.Cells.Item($row, 5).Value = $p
    $ws.Cells.Item($row, 6).Value = 111020
    $ws.Cells.Item($row, 7).Value = 2775509
    $row = $row + 1
}
